$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1. Rename header columns (row 1)
$ws.Range("A1").Value = "mx_state"
$ws.Range("B1").Value = "mx_municipality"
$ws.Range("C1").Value = "n_matriculas"
$ws.Range("D1").Value = "pct_matriculas"

# 2. Capitalize " de "/" el " -> " De "/" El " in specific municipality/state names
$ws.Range("B3").Value   = "Rincón De Romos"
$ws.Range("A20").Value  = "Ciudad De México"
$ws.Range("A30").Value  = "Estado De México"
$ws.Range("B38").Value  = "Acapulco De Juárez"
$ws.Range("B39").Value  = "Tlapa De Comonfort"
$ws.Range("B42").Value  = "Tulancingo De Bravo"
$ws.Range("B50").Value  = "Tamazula De Gordiano"
$ws.Range("B52").Value  = "Valle De Guadalupe"
$ws.Range("B53").Value  = "Zapotitlán De Vadillo"
$ws.Range("B77").Value  = "Heroica Ciudad De Huajuapan De León"
$ws.Range("B78").Value  = "Mariscala De Juárez"
$ws.Range("B79").Value  = "Oaxaca De Juárez"
$ws.Range("B80").Value  = "Putla Villa De Guerrero"
$ws.Range("B85").Value  = "Tlacolula De Matamoros"
$ws.Range("B88").Value  = "Huehuetlán El Grande"
$ws.Range("B94").Value  = "Amealco De Bonfil"
$ws.Range("B97").Value  = "Mexquitic De Carmona"
$ws.Range("B113").Value = "Poza Rica De Hidalgo"
$ws.Range("B120").Value = "Tlaltenango De Sánchez Román"

# 3. Tiny floating point precision fix on D54
$ws.Range("D54").Value = 0.09316770186335405

# 4. Remove trailing footnote/metadata rows 126-130
$ws.Range("A126:A130").EntireRow.Delete()
